$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (bold, centered, bordered) from H1 into the
# two new header cells, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new column I (I0) and column J (IF)
$dataI = @(9, 6, 6, 5, 5, 5, 7, 5, 6, 8, 5, 5)
$dataJ = @(9, 7, 7, 6, 6, 7, 8, 5, 6, 8, 5, 5)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
